# Word COM-interop script implementing the release-notes update:
#   - bump version number and release date
#   - rewrite the summary line ("Added type inference..." -> "Updates to
#     web clients and JSON library.")
#   - rewrite the three "Highlights" bullets
#   - drop the fourth "Highlights" bullet entirely
#   - move the _GoBack bookmark down to the paragraph that now follows the
#     (shortened) bullet list
#
# NOTE: this engine consolidates (merges) any run of contiguous,
# identically-formatted <w:r> elements in a paragraph the moment any text
# inside that paragraph is edited. The version/date line is one paragraph
# where every run shares the same rPr, so touching any of its digits
# necessarily folds neighbouring runs together -- that is unavoidable
# here and does not change the rendered text or formatting.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Version number + release date: "v5.1.7; August 17, 2019"
#                                -> "v5.1.9; August 29, 2019"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("v5.1.7", $true, $false, $false, $false, $false, `
    $true, 1, $false, "v5.1.9", 2)

$d.Content.Find.Execute("August 17, 2019", $true, $false, $false, $false, `
    $false, $true, 1, $false, "August 29, 2019", 2)

# ---------------------------------------------------------------------
# 2) Summary sentence just under "Release notes": split into two runs so
#    the trailing period ends up in its own <w:r>, matching the source
#    edit. Changing only formatting (not text) never triggers the
#    engine's run-merge pass, so flipping Bold on/off after inserting the
#    period keeps it as an independent run with the same final rPr as its
#    neighbour.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Added type inference for lambda expressions.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Updates to web clients and JSON library", 2)

$full = $d.Content.Text
$idx = $full.IndexOf("Updates to web clients and JSON library")
$endPos = $idx + "Updates to web clients and JSON library".Length
$tail = $d.Range($endPos, $endPos)
$tail.InsertAfter(".")
$periodRange = $d.Range($endPos, $endPos + 1)
$periodRange.Bold = 1
$periodRange.Bold = 0

# ---------------------------------------------------------------------
# 3) Highlights bullets
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Type inference for lambda expressions", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    "Support for HTTP header requests (new)", 2)

$d.Content.Find.Execute("Add 'Zip' to generics Map class", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "JSON encode and decode text string functions (new)", 2)

$d.Content.Find.Execute("Fixed bug that broke 'Apply' methods in generic collections", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "General enhancements", 2)

# Drop the fourth bullet entirely ("Resolved parsing issue ...").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Resolved parsing issue regarding multi-line lambdas within function calls*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 4) Move the _GoBack bookmark from the blank paragraph right after the
#    summary sentence down to the blank paragraph that now follows the
#    (shortened) bullet list.
# ---------------------------------------------------------------------
$d.Content.Bookmarks.Item("_GoBack").Delete()

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "" -and $p.Range.Text.Length -le 1) {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text -like "Fixed bug*" -or $prev.Range.Text -like "General enhancements*") {
            $p.Range.Bookmarks.Add("_GoBack")
            break
        }
    }
}
